$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95, pushing the existing rows 95:126 down to 96:127
$ws.Rows.Item(95).Insert()

# Make sure the new row's date cell (D95) keeps the same date number format
# used by the rest of the "Fecha" column (it is copied from the row above by
# Insert(), but set it explicitly as well to be safe).
$ws.Cells.Item(95,4).NumberFormat = $ws.Cells.Item(96,4).NumberFormat

# Populate the new row with the record that was added by this edit
$ws.Cells.Item(95,1).Value  = 9
$ws.Cells.Item(95,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(95,3).Value  = "Metropolitana"
$ws.Cells.Item(95,4).Value  = 44964
$ws.Cells.Item(95,5).Value  = 13
$ws.Cells.Item(95,6).Value  = "Fruta"
$ws.Cells.Item(95,7).Value  = 100101
$ws.Cells.Item(95,8).Value  = "Berries"
$ws.Cells.Item(95,9).Value  = 100101004
$ws.Cells.Item(95,10).Value = "Frambuesa"
$ws.Cells.Item(95,11).Value = "Sin especificar"
$ws.Cells.Item(95,12).Value = "Primera"
$ws.Cells.Item(95,13).Value = 350
$ws.Cells.Item(95,14).Value = 8000
$ws.Cells.Item(95,15).Value = 8000
$ws.Cells.Item(95,16).Value = 8000
$ws.Cells.Item(95,17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(95,18).Value = "Región de O'Higgins"
$ws.Cells.Item(95,19).Value = 4000
$ws.Cells.Item(95,20).Value = 2
